# Append a new data row (row 15) to the sheet, mirroring the pattern of
# the existing rows (A: date/time, B-M: numeric stats, N: "Bag" label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = 42620.885636574072
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = 18
$ws.Cells.Item($row, 3).Value = 59
$ws.Cells.Item($row, 4).Value = 35
$ws.Cells.Item($row, 5).Value = 59
$ws.Cells.Item($row, 6).Value = 14
$ws.Cells.Item($row, 7).Value = 14823
$ws.Cells.Item($row, 8).Value = 30550
$ws.Cells.Item($row, 9).Value = 3473
$ws.Cells.Item($row, 10).Value = 492
$ws.Cells.Item($row, 11).Value = 297
$ws.Cells.Item($row, 12).Value = 52
$ws.Cells.Item($row, 13).Value = 9
$ws.Cells.Item($row, 14).Value = "Bag"
